# Updated capital structure database
# Applies the refreshed data for the Finland / Financial Svcs. (Non-bank & Insurance)
# sheet: row 2 (aggregate), row 3 (Fellow Finance Oyj), row 4 / row 5 (company names
# swapped: row 4 becomes United Bankers Oyj, row 5 becomes Ferratum Oyj) along with
# all of their refreshed metric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2 - aggregate / industry row
# ---------------------------------------------------------------------------
$ws.Range("D2").Value = 0.1063
$ws.Range("E2").Value = 0.03885
$ws.Range("F2").Value = 0.03700000000000001
$ws.Range("G2").Value = 0.09886471144749291
$ws.Range("H2").Value = 0.008987701040681174
$ws.Range("I2").Value = 0.001182592242194891
$ws.Range("J2").Value = 0.0009910106997811686
$ws.Range("K2").Value = 15.06
$ws.Range("L2").Value = 0.07123935666982024
$ws.Range("M2").Value = 5.31
$ws.Range("N2").Value = 0.01798780487804878
$ws.Range("O2").Value = 0.352589641434263
$ws.Range("P2").Value = 5.05
$ws.Range("Q2").Value = 0.0171070460704607
$ws.Range("R2").Value = 0.3353253652058433
$ws.Range("S2").Value = 0.2599999999999998
$ws.Range("T2").Value = 0.04896421845574384
$ws.Range("U2").Value = 317.96
$ws.Range("V2").Value = 1.07710027100271
$ws.Range("W2").Value = 0.06659090909090909
$ws.Range("X2").Value = 0.02114035451499249
$ws.Range("Y2").Value = 0.0454505545759166
$ws.Range("Z2").Value = 0.7797486656855259
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0.02012841353003278
$ws.Range("AC2").Value = -0.01784852803162772
$ws.Range("AD2").Value = 227.53
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 227.53
$ws.Range("AG2").Value = -90.42999999999998
$ws.Range("AH2").Value = 0.4352725116216785
$ws.Range("AI2").Value = 0.523743756186267
$ws.Range("AJ2").Value = -0.4416174244274063
$ws.Range("AK2").Value = -0.7764231132480464
$ws.Range("AL2").Value = 1.66
$ws.Range("AM2").Value = 1.101
$ws.Range("AN2").Value = 651.948424068768
$ws.Range("AO2").Value = 0.1506024096385542
$ws.Range("AP2").Value = -259.1117478510028
$ws.Range("AQ2").Value = 0.2270663033605813

# ---------------------------------------------------------------------------
# Row 3 - Fellow Finance Oyj (HLSE:FELLOW)
# ---------------------------------------------------------------------------
$ws.Range("F3").ClearContents()
$ws.Range("G3").Value = 0.1319444444444444
$ws.Range("H3").Value = 0.1319444444444444
$ws.Range("I3").Value = 0.01736111111111111
$ws.Range("J3").Value = 0.01736111111111111
$ws.Range("K3").Value = -1.37
$ws.Range("L3").Value = -0.0951388888888889
$ws.Range("M3").Value = -0
$ws.Range("N3").Value = -0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = -0
$ws.Range("Q3").Value = -0
$ws.Range("R3").Value = 0
$ws.Range("T3").ClearContents()
$ws.Range("U3").Value = 4.46
$ws.Range("V3").Value = 0.1755905511811024
$ws.Range("W3").Value = -0.08203592814371259
$ws.Range("X3").Value = 0.02114035451499249
$ws.Range("Y3").Value = -0.1031762826587051
$ws.Range("Z3").Value = 1.611820013431834
$ws.Range("AA3").Value = 0.02798298634430267
$ws.Range("AB3").Value = 0.02012841353003278
$ws.Range("AC3").Value = 0.007854572814269885
$ws.Range("AD3").Value = 13.1
$ws.Range("AF3").Value = 13.1
$ws.Range("AG3").Value = 8.640000000000001
$ws.Range("AH3").Value = 0.3402597402597403
$ws.Range("AI3").Value = 0.458041958041958
$ws.Range("AJ3").Value = 0.2538190364277321
$ws.Range("AK3").Value = 0.3579121789560895
$ws.Range("AL3").Value = 1.66
$ws.Range("AM3").Value = 1.636
$ws.Range("AN3").Value = 37.53581661891118
$ws.Range("AO3").Value = 0.1506024096385542
$ws.Range("AP3").Value = 24.75644699140402
$ws.Range("AQ3").Value = 0.1528117359413203

# ---------------------------------------------------------------------------
# Row 4 - now United Bankers Oyj (HLSE:UNIAV) (previously Ferratum Oyj)
# ---------------------------------------------------------------------------
$ws.Range("B4").Value = "United Bankers Oyj (HLSE:UNIAV)"
$ws.Range("D4").Value = 0.05860000000000001
$ws.Range("E4").Value = 0.12
$ws.Range("G4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 7.64
$ws.Range("L4").Value = 0.1919597989949749
$ws.Range("M4").Value = 5.31
$ws.Range("N4").Value = 0.03898678414096916
$ws.Range("O4").Value = 0.6950261780104712
$ws.Range("P4").Value = 5.05
$ws.Range("Q4").Value = 0.03707782672540382
$ws.Range("R4").Value = 0.6609947643979057
$ws.Range("S4").Value = 0.2599999999999998
$ws.Range("T4").Value = 0.04896421845574384
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0.2301204819277108
$ws.Range("X4").Value = 0.01784099133054493
$ws.Range("Y4").Value = 0.2122794905971659
$ws.Range("Z4").Value = 0.9832258702043032
$ws.Range("AA4").Value = 0
$ws.Range("AB4").Value = 0.01784852803162772
$ws.Range("AC4").Value = -0.01784852803162772
$ws.Range("AD4").Value = 3.23
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 3.23
$ws.Range("AG4").Value = 3.23
$ws.Range("AH4").Value = 0.02316574625259988
$ws.Range("AI4").Value = 0.07221104404202996
$ws.Range("AJ4").Value = 0.02316574625259988
$ws.Range("AK4").Value = 0.07221104404202996
$ws.Range("AM4").Value = 0
$ws.Range("AN4").ClearContents()
$ws.Range("AP4").ClearContents()
$ws.Range("AQ4").ClearContents()

# ---------------------------------------------------------------------------
# Row 5 - now Ferratum Oyj (XTRA:FRU) (previously United Bankers Oyj)
# ---------------------------------------------------------------------------
$ws.Range("B5").Value = "Ferratum Oyj (XTRA:FRU)"
$ws.Range("D5").Value = 0.154
$ws.Range("E5").Value = -0.0423
$ws.Range("F5").Value = 0.03700000000000001
$ws.Range("G5").Value = 0.1208651399491094
$ws.Range("K5").Value = 8.789999999999999
$ws.Range("L5").Value = 0.05591603053435115
$ws.Range("M5").Value = -0
$ws.Range("N5").Value = -0
$ws.Range("O5").Value = -0
$ws.Range("P5").Value = -0
$ws.Range("Q5").Value = -0
$ws.Range("R5").Value = -0
$ws.Range("S5").Value = 0
$ws.Range("T5").ClearContents()
$ws.Range("U5").Value = 313.5
$ws.Range("V5").Value = 2.346556886227545
$ws.Range("W5").Value = 0.06659090909090909
$ws.Range("X5").Value = 0.02828239653974361
$ws.Range("Y5").Value = 0.03830851255116548
$ws.Range("Z5").Value = 0.7090663058186739
$ws.Range("AA5").Value = 0
$ws.Range("AB5").Value = 0.02532016521023656
$ws.Range("AC5").Value = -0.02532016521023656
$ws.Range("AD5").Value = 211.2
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 211.2
$ws.Range("AG5").Value = -102.3
$ws.Range("AH5").Value = 0.6125290023201857
$ws.Range("AI5").Value = 0.5848795347549155
$ws.Range("AJ5").Value = -3.268370607028756
$ws.Range("AK5").Value = -2.149159663865547
$ws.Range("AM5").Value = -0.535
$ws.Range("AQ5").Value = -0
